$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.472.68'
$ws.Range("E2").Value = '  -7.73%  '

# Row 3
$ws.Range("D3").Value = '3.672.32'
$ws.Range("E3").Value = '  -7.82%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''577.50'
$ws.Range("E5").Value = '  -5.67%  '

# Row 6
$ws.Range("D6").Value = '''170.70'
$ws.Range("E6").Value = '  -1.37%  '

# Row 7
$ws.Range("D7").Value = '3.658.13'
$ws.Range("E7").Value = '  -7.97%  '

# Row 8
$ws.Range("D8").Value = '''0.621'
$ws.Range("E8").Value = '  -10.39%  '

# Row 9
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.15%  '

# Row 11
$ws.Range("E11").Value = '  -13.61%  '

# Row 12
$ws.Range("D12").Value = '''51.26'
$ws.Range("E12").Value = '  -10.57%  '

# Row 13
$ws.Range("D13").Value = '''0.0000289'
$ws.Range("E13").Value = '  -14.38%  '

# Row 14
$ws.Range("D14").Value = '''10.40'
$ws.Range("E14").Value = '  -11.93%  '

# Row 15
$ws.Range("D15").Value = '4.265.01'
$ws.Range("E15").Value = '  -7.66%  '

# Row 16
$ws.Range("D16").Value = '3.704.35'
$ws.Range("E16").Value = '  -6.98%  '

# Row 17
$ws.Range("D17").Value = '''19.28'
$ws.Range("E17").Value = '  -10.41%  '

# Row 18
$ws.Range("E18").Value = '  -3.62%  '

# Row 19
$ws.Range("D19").Value = '''12.79'
$ws.Range("E19").Value = '  -10.85%  '

# Row 20
$ws.Range("D20").Value = '''1.11'
$ws.Range("E20").Value = '  -11.21%  '

# Row 21
$ws.Range("D21").Value = '67.404.78'
$ws.Range("E21").Value = '  -7.76%  '

# Row 22
$ws.Range("D22").Value = '''404.05'
$ws.Range("E22").Value = '  -11.76%  '

# Row 23
$ws.Range("D23").Value = '''4.46'
$ws.Range("E23").Value = '  -8.53%  '

# Row 24
$ws.Range("D24").Value = '''87.29'
$ws.Range("E24").Value = '  -10.49%  '

# Row 25
$ws.Range("D25").Value = '''3.02'
$ws.Range("E25").Value = '  -10.92%  '

# Row 26
$ws.Range("D26").Value = '''12.69'
$ws.Range("E26").Value = '  -11.63%  '

# Row 27
$ws.Range("D27").Value = '''10.57'
$ws.Range("E27").Value = '  -7.43%  '

# Row 28
$ws.Range("E28").Value = '  +1.74%  '

# Row 29
$ws.Range("D29").Value = '''3.73'
$ws.Range("E29").Value = '  -12.74%  '

# Row 30
$ws.Range("D30").Value = '''9.39'
$ws.Range("E30").Value = '  -11.94%  '

# Row 31
$ws.Range("D31").Value = '''32.48'
$ws.Range("E31").Value = '  -10.68%  '

# Row 32
$ws.Range("D32").Value = '''7.45'
$ws.Range("E32").Value = '  -6.59%  '

# Row 33
$ws.Range("D33").Value = '''12.39'
$ws.Range("E33").Value = '  -12.16%  '

# Row 34
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '''64.97'
$ws.Range("E34").Value = '  -6.64%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.116'
$ws.Range("E35").Value = '  -10.80%  '

# Row 36
$ws.Range("D36").Value = '''42.91'
$ws.Range("E36").Value = '  -13.73%  '

# Row 37
$ws.Range("D37").Value = '''593.42'
$ws.Range("E37").Value = '  -6.42%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0880'
$ws.Range("E38").Value = '  -15.11%  '

# Row 39
$ws.Range("E39").Value = '  -0.06%  '

# Row 40
$ws.Range("D40").Value = '''0.394'
$ws.Range("E40").Value = '  -8.94%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$ws.Range("E42").Value = '  -10.60%  '

# Row 43
$ws.Range("D43").Value = '''2.97'
$ws.Range("E43").Value = '  -13.90%  '

# Row 44
$ws.Range("D44").Value = '''2.65'
$ws.Range("E44").Value = '  -0.91%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0434'
$ws.Range("E45").Value = '  -11.03%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '''2.88'
$ws.Range("E46").Value = '  -12.17%  '

# Row 47
$ws.Range("D47").Value = '''9.12'
$ws.Range("E47").Value = '  -13.83%  '

# Row 48
$ws.Range("D48").Value = '2.783.60'
$ws.Range("E48").Value = '  -1.13%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.133'
$ws.Range("E49").Value = '  -11.46%  '

# Row 50
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '''2.69'
$ws.Range("E50").Value = '  -10.43%  '

# Row 51
$ws.Range("D51").Value = '''3.16'
$ws.Range("E51").Value = '  -7.63%  '

# Strip the quote-prefix formatting flag introduced by forcing text above,
# restoring default (unstyled) presentation while keeping literal text values.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
